$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Main")

# Update the third-party vote reallocation inputs:
#   pct_stein_to_clinton (D3): 0.8 -> 0.9
#   pct_johnson_to_clinton (E3): 0.64 -> 0.6
$ws.Range("D3").Value = 0.9
$ws.Range("E3").Value = 0.6

# Move/restore the active selection on the bottom-right frozen pane to D3
$ws.Range("D3").Select()
